# Renderer017-TemplateFormula/template.xlsx update:
# "feat: update foreach/endrow/endloop with new behaviour"
#
# Adds a new template-directive column (E) demonstrating the new
# "#! END_ROW true" / "#! END_ROW" behaviour to both worksheets, and
# switches the active sheet/selection to reflect where the author was
# last working ("Incorrect Formula" sheet, cell D13) while the
# "Correct Formula" sheet keeps a parked selection at F14.

$wb = $excel.ActiveWorkbook

$wsCorrect = $wb.Worksheets.Item("Correct Formula")
$wsIncorrect = $wb.Worksheets.Item("Incorrect Formula")

# --- New column E: "#! END_ROW true" header directive + "#! END_ROW" row ---
$wsCorrect.Range("E1").Value = "#! END_ROW true"
$wsCorrect.Range("E2").Value = "#! END_ROW"

$wsIncorrect.Range("E1").Value = "#! END_ROW true"
$wsIncorrect.Range("E2").Value = "#! END_ROW"

# --- Selections left on each sheet ---
$wsCorrect.Range("F14").Select()
$wsIncorrect.Range("D13").Select()

# --- Make "Incorrect Formula" the active (visible) tab ---
$wsIncorrect.Activate()
